$wb = $excel.ActiveWorkbook

# Hunk 0: ALC!row15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 216.81
$ws.Range("I15").Value = 216.81
$ws.Range("K15").Value = 650.4300000000001
$ws.Range("M15").Value = -481.4300000000001

# Hunk 1: ALC!row76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 71431336
$ws.Range("I76").Value = 83336000
$ws.Range("J76").Value = 3333.3333
$ws.Range("K76").Value = 83336000
$ws.Range("L76").Value = 3333.3333
$ws.Range("M76").Value = -83335685
$ws.Range("N76").Value = -3963.3333

# Hunk 2: ALC!row79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 71431336
$ws.Range("I79").Value = 83336000
$ws.Range("J79").Value = 3333.3333
$ws.Range("K79").Value = 83336000
$ws.Range("L79").Value = 3333.3333
$ws.Range("M79").Value = -83334908
$ws.Range("N79").Value = -5517.3333

# Hunk 3: ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4776.4287
$ws.Range("I116").Value = 4655.8335
$ws.Range("J116").Value = 5500
$ws.Range("K116").Value = 4655.8335
$ws.Range("L116").Value = 5500
$ws.Range("M116").Value = -1213.8335
$ws.Range("N116").Value = -12384

# Hunk 4: ARM!row29
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 12000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 12000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 12000
$ws.Range("N29").Value = -12616
$ws.Range("M29").ClearContents()

# Hunk 5: ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1763.4706
$ws.Range("I132").Value = 1356.3478
$ws.Range("K132").Value = 4069.0434
$ws.Range("M132").Value = -1539.0434

# Hunk 6: CRP!row16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1552.95
$ws.Range("I16").Value = 1233.3
$ws.Range("J16").Value = 1872.6
$ws.Range("K16").Value = 1233.3
$ws.Range("L16").Value = 1872.6
$ws.Range("M16").Value = -946.3
$ws.Range("N16").Value = -2446.6

# Hunk 7: CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37479.863
$ws.Range("I31").Value = 3100.0833
$ws.Range("J31").Value = 202502.8
$ws.Range("K31").Value = 3100.0833
$ws.Range("L31").Value = 202502.8
$ws.Range("M31").Value = -2805.0833
$ws.Range("N31").Value = -203092.8

# Hunk 8: CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 37479.863
$ws.Range("I34").Value = 3100.0833
$ws.Range("J34").Value = 202502.8
$ws.Range("K34").Value = 3100.0833
$ws.Range("L34").Value = 202502.8
$ws.Range("M34").Value = -2898.0833
$ws.Range("N34").Value = -202906.8

# Hunk 9: CRP!row107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 324.55554
$ws.Range("I107").Value = 342.1
$ws.Range("J107").Value = 314.2353
$ws.Range("K107").Value = 342.1
$ws.Range("L107").Value = 314.2353
$ws.Range("M107").Value = 1577.9
$ws.Range("N107").Value = -4154.2353

# Hunk 10: CRP!row113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1552.95
$ws.Range("I113").Value = 1233.3
$ws.Range("J113").Value = 1872.6
$ws.Range("K113").Value = 1233.3
$ws.Range("L113").Value = 1872.6
$ws.Range("M113").Value = 936.7
$ws.Range("N113").Value = -6212.6

# Hunk 11: CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1494.4166
$ws.Range("I132").Value = 1215.5
$ws.Range("J132").Value = 2052.25
$ws.Range("K132").Value = 3646.5
$ws.Range("L132").Value = 6156.75
$ws.Range("M132").Value = -1116.5
$ws.Range("N132").Value = -11216.75

# Hunk 12: CUL!row3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3927.75
$ws.Range("I3").Value = 929.8125
$ws.Range("K3").Value = 2789.4375
$ws.Range("M3").Value = -2677.4375

# Hunk 13: CUL!row23
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 144.4
$ws.Range("I23").Value = 80
$ws.Range("J23").Value = 176.6
$ws.Range("K23").Value = 240
$ws.Range("L23").Value = 529.8
$ws.Range("M23").Value = -5
$ws.Range("N23").Value = -999.8

# Hunk 14: CUL!row68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 13355649
$ws.Range("I68").Value = 31145116
$ws.Range("J68").Value = 13550
$ws.Range("K68").Value = 93435348
$ws.Range("L68").Value = 40650
$ws.Range("M68").Value = -93434537
$ws.Range("N68").Value = -42272

# Hunk 15: CUL!row71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 13355649
$ws.Range("I71").Value = 31145116
$ws.Range("J71").Value = 13550
$ws.Range("K71").Value = 280306044
$ws.Range("L71").Value = 121950
$ws.Range("M71").Value = -280301988
$ws.Range("N71").Value = -130062

# Hunk 16: CUL!row80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 9185
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 9782.857
$ws.Range("K80").Value = 15000
$ws.Range("L80").Value = 29348.571
$ws.Range("M80").Value = -14064
$ws.Range("N80").Value = -31220.571

# Hunk 17: CUL!row83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 9185
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 9782.857
$ws.Range("K83").Value = 45000
$ws.Range("L83").Value = 88045.713
$ws.Range("M83").Value = -40320
$ws.Range("N83").Value = -97405.713

# Hunk 18: GSM!row10
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 5002501.5
$ws.Range("I10").Value = 5002501.5
$ws.Range("K10").Value = 5002501.5
$ws.Range("M10").Value = -5002332.5

# Hunk 19: GSM!row18
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 2049203.4
$ws.Range("J18").Value = 78670.664
$ws.Range("L18").Value = 78670.664
$ws.Range("N18").Value = -79256.664

# Hunk 20: GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4395.45
$ws.Range("I70").Value = 3991.6667
$ws.Range("K70").Value = 3991.6667
$ws.Range("M70").Value = -3721.6667

# Hunk 21: GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4395.45
$ws.Range("I73").Value = 3991.6667
$ws.Range("K73").Value = 3991.6667
$ws.Range("M73").Value = -3055.6667

# Hunk 22: GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1429.4286
$ws.Range("I102").Value = 1401.2
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1401.2
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 220.8
$ws.Range("N102").Value = -4744

# Hunk 23: GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2851.077
$ws.Range("I126").Value = 2938.6667
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 8816.000100000001
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -6346.000100000001
$ws.Range("N126").Value = -10340

# Hunk 24: LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3277.8333
$ws.Range("I132").Value = 2261.8462
$ws.Range("K132").Value = 6785.5386
$ws.Range("M132").Value = -4255.5386

# Hunk 25: WVR!row14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 35000
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# Hunk 26: WVR!row25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

# Hunk 27: WVR!row107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 695.9474
$ws.Range("I107").Value = 651.4286
$ws.Range("J107").Value = 820.6
$ws.Range("K107").Value = 1954.2858
$ws.Range("L107").Value = 2461.8
$ws.Range("M107").Value = -34.28579999999988
$ws.Range("N107").Value = -6301.8

# Hunk 28: WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 952.2162
$ws.Range("I132").Value = 964.7778
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 2894.3334
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -364.3334
$ws.Range("N132").Value = -6560
